$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values of D2/D3, K2/K3, L2/L3, M2/M3, O2/O3, P2/P3
$cols = @("D", "K", "L", "M", "O", "P")

foreach ($col in $cols) {
    $addr2 = "{0}2" -f $col
    $addr3 = "{0}3" -f $col
    $v2 = $ws.Range($addr2).Value2
    $v3 = $ws.Range($addr3).Value2
    $ws.Range($addr2).Value2 = $v3
    $ws.Range($addr3).Value2 = $v2
}
